$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing phone numbers in column A (drop the leading "55" country code) ---
# Set the numeric value first, THEN apply the Text format, so the stored cell
# keeps its numeric type (matches target: <c r="A2" s="1"><v>61998232332</v></c>)
$ws.Range("A2").Value2 = 61998232332
$ws.Range("A3").Value2 = 61998232332
$ws.Range("A4").Value2 = 61998232332

# --- Apply the "Text" number format to the whole of column A (entire column) ---
# Done now, while only rows 1-4 exist, so it doesn't create phantom cells in
# the still-to-be-added rows 5-7.
$ws.Range("A1:A1048576").NumberFormat = "@"

# --- Resize column A (mirrors selecting the column header + Format > Column Width) ---
$ws.Columns.Item(1).ColumnWidth = 60.17
$ws.Columns.Item(1).Select() | Out-Null

# --- New row 5: blank contact, same message as rows 3/4 ---
$ws.Range("B5").Value2 = "Boa tarde, fique aqui com esta mensagem:"

# --- New row 6: contact typed as literal text (format applied before typing) ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value2 = "61998232332"
$ws.Range("B6").Value2 = "Boa tarde, fique aqui com esta mensagem:"

# --- New row 7: an empty, underlined, text-formatted placeholder cell ---
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Font.Underline = $true

# --- Page setup: A4 paper, portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
